$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two changed data values (Log K column) ---
$ws.Range("B3").Value = -2.9
$ws.Range("B4").Value = 4.5999999999999996

# --- Update the "best fit" column widths for columns A and C ---
# (Target stored widths 48.375 / 12.375 correspond to ColumnWidth 47.5 / 11.5
#  once Excel re-quantizes to its internal Normal-style character grid.)
$ws.Columns.Item(1).ColumnWidth = 47.5
$ws.Columns.Item(3).ColumnWidth = 11.5

# --- Update the view: scroll the window down so row 13 is at the top,
#     and move the active selection to B20 ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B20").Select()
